# Applies the BOQ/bill updates described by the diff:
# rows 8-12 get new quantities/descriptions/rates/amounts, and the
# Grand Total / Net Payable rows (14 & 16) are refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 ---------------------------------------------------------------
$ws.Range("C8").Value = 17
$ws.Range("D8").Value = "'3"
$ws.Range("E8").Value = "Medium point (up to 6 mtr.)"
$ws.Range("F8").Value = 472
$ws.Range("G8").Value = "'8024.00"

# --- Row 9 -----------------------------------------------------------------
$ws.Range("A9").Value = "Each"
$ws.Range("C9").Value = 66
$ws.Range("D9").Value = "'3.0"
$ws.Range("E9").Value = 'P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F9").Value = 23
$ws.Range("G9").Value = "'1518.00"

# --- Row 10 ------------------------------------------------------------
$ws.Range("A10").Value = "R. mtr."
$ws.Range("C10").Value = 64
$ws.Range("D10").Value = "'16"
$ws.Range("E10").Value = "20 mm"
$ws.Range("F10").Value = 40
$ws.Range("G10").Value = "'2560.00"

# --- Row 11 ------------------------------------------------------------
$ws.Range("A11").Value = "Each"
$ws.Range("C11").Value = 38
$ws.Range("D11").Value = "'27"
$ws.Range("E11").Value = "1170mm(+/-10%) LED batten with min. lumen output 2200 lm"
$ws.Range("F11").Value = 492
$ws.Range("G11").Value = "'18696.00"

# --- Row 12 ------------------------------------------------------------
$ws.Range("C12").Value = 24
$ws.Range("D12").Value = "'30"
$ws.Range("E12").Value = " 6 A to 32 A rating"
$ws.Range("F12").Value = 187
$ws.Range("G12").Value = "'4488.00"

# --- Totals (Grand Total / Net Payable) -----------------------------------
$ws.Range("G14").Value = "'35286.00"
$ws.Range("H14").Value = "'35286.00"

$ws.Range("G16").Value = "'35286.00"
$ws.Range("H16").Value = "'35286.00"
